$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -------------------------------------------------------------------------
# The old row 19 ("Are you sure you wanna go out?" / "Are you sure you wanna
# go out?" confirmation pair) duplicated the content already present in row
# 5, so it is removed entirely. Deleting the row also shifts every
# subsequent row up by one and drops the now-unused "problem getting video"
# strings from the shared-string table.
# -------------------------------------------------------------------------
$ws.Rows("19:19").Delete()

# After the delete:
#   row 18 -> Redes sociales / Social Media            (unchanged)
#   row 19 -> Hay una nueva version disponible / ...    (was row 20)
#   row 20 -> Quieres actualizar... / Do you wanna ...  (was row 21)
#   row 21 -> Mas tarde / Later                         (was row 22)
#   row 22 -> now the blank trailing row (was row 23)

# -------------------------------------------------------------------------
# New feature strings: "Check for updates automatically" checkbox option,
# plus the "already up to date" notice dialog.
# -------------------------------------------------------------------------

# Row 22 used to be the empty trailing row - it now holds the new checkbox
# label pair, keeping the same (center-aligned, no wrap) style as row 21.
$ws.Range("A22").Value = "Comprobar actualizaciones automáticamente"
$ws.Range("B22").Value = "Check for updates automatically"
$ws.Range("A22:B22").HorizontalAlignment = -4108  # xlCenter
$ws.Range("A22:B22").WrapText = $false

# Rows 23 and 24 are brand new, appended at the bottom of the table.
$ws.Range("A23").Value = "Versión Actualizada"
$ws.Range("B23").Value = "Updated Version"
$ws.Range("A23:B23").HorizontalAlignment = -4108  # xlCenter
$ws.Range("A23:B23").WrapText = $false

$ws.Range("A24").Value = "¡Felicidades! Estás utilizando la última versión de la aplicación. " + [char]10 + "No es necesario realizar ninguna actualización en este momento."
$ws.Range("B24").Value = "Congratulations! You are using the latest version of the app. " + [char]10 + "No update is required at this time."
$ws.Range("A24:B24").HorizontalAlignment = -4108  # xlCenter
$ws.Range("A24:B24").WrapText = $true
$ws.Rows("24:24").RowHeight = 45

# -------------------------------------------------------------------------
# Update the view state: scroll so row 10 is at the top and select C17
# (cosmetic - matches the author's last on-screen position).
# -------------------------------------------------------------------------
$ws.Range("C17").Select()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
